$d = $word.ActiveDocument

# The edit: place the cursor at the end of the document (end of the
# "Nuevo cambio" paragraph), press Enter to start a new paragraph, and
# type the new sentence there. Word keeps its "_GoBack" bookmark glued to
# the most recent edit location, so we recreate it at the end of the text
# we just typed.

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$rng = $lastPara.Range
$rng.Collapse(0)              # wdCollapseEnd
$rng.InsertParagraphAfter()
$rng.Collapse(0)

$newText = "Esto se escribió la segunda vez que se abrió git"

# Type the new text plus one throw-away trailing character. A collapsed
# bookmark placed exactly at the end of the last paragraph of the story
# is mishandled by this host, so we park the bookmark just before that
# trailing character (a perfectly ordinary, non-boundary position) and
# then delete the character, leaving the bookmark collapsed right after
# the real text - exactly where Word would leave "_GoBack".
$rng.InsertAfter($newText + "#")

$markerPos = $rng.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($markerPos, $markerPos))
$d.Range($markerPos, $markerPos + 1).Delete()
